$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3806534858468353
$ws.Range("D2").Value = 0.03833107441264616
$ws.Range("E2").Value = 0.1531497260172436
$ws.Range("F2").Value = 0.7282730210210957
$ws.Range("G2").Value = 0.5703386385883817
$ws.Range("H2").Value = 0.7097870972761982
$ws.Range("K2").Value = 1.267548457364114
$ws.Range("L2").Value = 0.1301955329974689
$ws.Range("M2").Value = 0.3998360253669802
$ws.Range("O2").Value = 2.536875369564072
$ws.Range("C3").Value = 0.3779189058160028
$ws.Range("D3").Value = 0.03576922354120882
$ws.Range("E3").Value = 0.1535404883569207
$ws.Range("F3").Value = 0.7335329962084316
$ws.Range("G3").Value = 0.5772903253341894
$ws.Range("H3").Value = 0.7179173613832361
$ws.Range("K3").Value = 1.114592041293292
$ws.Range("L3").Value = 0.1314559830054129
$ws.Range("M3").Value = 0.3673673862284659
$ws.Range("O3").Value = 2.568326242273088
$ws.Range("C4").Value = 0.3764499214083088
$ws.Range("D4").Value = 0.03418704470035294
$ws.Range("E4").Value = 0.1538674515135199
$ws.Range("F4").Value = 0.7373460182820253
$ws.Range("G4").Value = 0.5820885879117554
$ws.Range("H4").Value = 0.7233174499215238
$ws.Range("K4").Value = 1.020366150735498
$ws.Range("L4").Value = 0.1323030877200075
$ws.Range("M4").Value = 0.3474534518894217
$ws.Range("O4").Value = 2.589604354527125
$ws.Range("C5").Value = 0.3759041718964227
$ws.Range("D5").Value = 0.0335400260097316
$ws.Range("E5").Value = 0.1540225798112367
$ws.Range("F5").Value = 0.7390463430446914
$ws.Range("G5").Value = 0.5841768568123697
$ws.Range("H5").Value = 0.7256205970893035
$ws.Range("K5").Value = 0.9818930059653042
$ws.Range("L5").Value = 0.1326666856159981
$ws.Range("M5").Value = 0.339344405732021
$ws.Range("O5").Value = 2.598769254270039
$ws.Range("C6").Value = 0.3758167455932835
$ws.Range("D6").Value = 0.0334324533555872
$ws.Range("E6").Value = 0.1540496608241142
$ws.Range("F6").Value = 0.7393375220538587
$ws.Range("G6").Value = 0.584531632305108
$ws.Range("H6").Value = 0.7260092263777054
$ws.Range("K6").Value = 0.9755000978059059
$ws.Range("L6").Value = 0.1327281716382416
$ws.Range("M6").Value = 0.3379982854072878
$ws.Range("O6").Value = 2.600320886700487
$ws.Range("C7").Value = 0.3764423470867371
$ws.Range("D7").Value = 0.03417832790499631
$ws.Range("E7").Value = 0.153869454999672
$ws.Range("F7").Value = 0.7373683566115332
$ws.Range("G7").Value = 0.5821162131878239
$ws.Range("H7").Value = 0.7233480957410308
$ws.Range("K7").Value = 1.01984758994962
$ws.Range("L7").Value = 0.132307916853339
$ws.Range("M7").Value = 0.3473440652848367
$ws.Range("O7").Value = 2.589725956816253
$ws.Range("C8").Value = 0.3796670345671487
$ws.Range("D8").Value = 0.03744968124996717
$ws.Range("E8").Value = 0.153266398470258
$ws.Range("F8").Value = 0.7299654554340478
$ws.Range("G8").Value = 0.5726254138679252
$ws.Range("H8").Value = 0.7125056893517652
$ws.Range("K8").Value = 1.21487495708277
$ws.Range("L8").Value = 0.1306149429488208
$ws.Range("M8").Value = 0.3886365890116608
$ws.Range("O8").Value = 2.547310893380114
$ws.Range("C9").Value = 0.3876558926654354
$ws.Range("D9").Value = 0.04379025648996304
$ws.Range("E9").Value = 0.1527744153348358
$ws.Range("F9").Value = 0.720087097871108
$ws.Range("G9").Value = 0.5582325590836916
$ws.Range("H9").Value = 0.6944836758019264
$ws.Range("K9").Value = 1.594766461461745
$ws.Range("L9").Value = 0.12787602278555
$ws.Range("M9").Value = 0.4697662634768065
$ws.Range("O9").Value = 2.47977791360627
$ws.Range("C10").Value = 0.394539768653118
$ws.Range("D10").Value = 0.04840152929165242
$ws.Range("E10").Value = 0.1528342619863707
$ws.Range("F10").Value = 0.7156717936794408
$ws.Range("G10").Value = 0.5502495061860486
$ws.Range("H10").Value = 0.6832207857252541
$ws.Range("K10").Value = 1.872212305618973
$ws.Range("L10").Value = 0.1262183985588443
$ws.Range("M10").Value = 0.5294475564014647
$ws.Range("O10").Value = 2.439745507808198
$ws.Range("C11").Value = 0.3978916319120174
$ws.Range("D11").Value = 0.05048874950612969
$ws.Range("E11").Value = 0.1529530501432035
$ws.Range("F11").Value = 0.7142834789312502
$ws.Range("G11").Value = 0.5471847145297204
$ws.Range("H11").Value = 0.6785270831532415
$ws.Range("K11").Value = 1.998049530886021
$ws.Range("L11").Value = 0.1255414186170185
$ws.Range("M11").Value = 0.5566108327104331
$ws.Range("O11").Value = 2.423625099492682
$ws.Range("C12").Value = 0.3991925500854734
$ws.Range("D12").Value = 0.05127758281317085
$ws.Range("E12").Value = 0.1530112011878195
$ws.Range("F12").Value = 0.7138471879986028
$ws.Range("G12").Value = 0.5461059932516861
$ws.Range("H12").Value = 0.6768115677033677
$ws.Range("K12").Value = 2.045644786425783
$ws.Range("L12").Value = 0.1252961573325777
$ws.Range("M12").Value = 0.566898430915785
$ws.Range("O12").Value = 2.417822181400624
$ws.Range("C13").Value = 0.398910967608856
$ws.Range("D13").Value = 0.05110776318686305
$ws.Range("E13").Value = 0.1529980916171034
$ws.Range("F13").Value = 0.7139371695790544
$ws.Range("G13").Value = 0.5463346693395721
$ws.Range("H13").Value = 0.6771782809465208
$ws.Range("K13").Value = 2.035396856040563
$ws.Range("L13").Value = 0.1253484850341771
$ws.Range("M13").Value = 0.5646827566032329
$ws.Range("O13").Value = 2.419058518464112
$ws.Range("C14").Value = 0.3979980253799056
$ws.Range("D14").Value = 0.05055367865783467
$ws.Range("E14").Value = 0.1529575703371364
$ws.Range("F14").Value = 0.7142457912291249
$ws.Range("G14").Value = 0.5470943250968645
$ws.Range("H14").Value = 0.6783847058997026
$ws.Range("K14").Value = 2.001966370616572
$ws.Range("L14").Value = 0.1255210183674755
$ws.Range("M14").Value = 0.5574571735918283
$ws.Range("O14").Value = 2.42314164158222
$ws.Range("C15").Value = 0.3974429410299081
$ws.Range("D15").Value = 0.05021408245150383
$ws.Range("E15").Value = 0.1529344648693538
$ws.Range("F15").Value = 0.7144464850664178
$ws.Range("G15").Value = 0.5475703057418286
$ws.Range("H15").Value = 0.6791317379256014
$ws.Range("K15").Value = 1.981481784704329
$ws.Range("L15").Value = 0.125628145443816
$ws.Range("M15").Value = 0.5530314685404676
$ws.Range("O15").Value = 2.425681967652821
$ws.Range("C16").Value = 0.3943251466315019
$ws.Range("D16").Value = 0.04826490990051013
$ws.Range("E16").Value = 0.1528283415514728
$ws.Range("F16").Value = 0.7157750249438664
$ws.Range("G16").Value = 0.5504612319351168
$ws.Range("H16").Value = 0.6835361863446252
$ws.Range("K16").Value = 1.863980779857172
$ws.Range("L16").Value = 0.1262641929023331
$ws.Range("M16").Value = 0.5276726059698262
$ws.Range("O16").Value = 2.44084115955468
$ws.Range("C17").Value = 0.3924688857008505
$ws.Range("D17").Value = 0.04706644097718282
$ws.Range("E17").Value = 0.1527866912812996
$ws.Range("F17").Value = 0.716749072032286
$ws.Range("G17").Value = 0.5523801283309524
$ws.Range("H17").Value = 0.6863483206322201
$ws.Range("K17").Value = 1.791799891918856
$ws.Range("L17").Value = 0.1266741384708432
$ws.Range("M17").Value = 0.5121189456344979
$ws.Range("O17").Value = 2.450676920644298
$ws.Range("C18").Value = 0.3914219570054627
$ws.Range("D18").Value = 0.04637613050296352
$ws.Range("E18").Value = 0.1527713534446882
$ws.Range("F18").Value = 0.7173676825129363
$ws.Range("G18").Value = 0.5535371500004587
$ws.Range("H18").Value = 0.688006240118554
$ws.Range("K18").Value = 1.750248287817897
$ws.Range("L18").Value = 0.1269171828853253
$ws.Range("M18").Value = 0.5031742379757276
$ws.Range("O18").Value = 2.45653090680274
$ws.Range("C19").Value = 0.3910710490604856
$ws.Range("D19").Value = 0.04614223583984511
$ws.Range("E19").Value = 0.1527676404614482
$ws.Range("F19").Value = 0.7175871500165556
$ws.Range("G19").Value = 0.553938045205733
$ws.Range("H19").Value = 0.6885745288750158
$ws.Range("K19").Value = 1.736173690150906
$ws.Range("L19").Value = 0.1270007191586373
$ws.Range("M19").Value = 0.5001459617909063
$ws.Range("O19").Value = 2.458546724260614
$ws.Range("C20").Value = 0.3926643413074373
$ws.Range("D20").Value = 0.04719412209555429
$ws.Range("E20").Value = 0.1527902330489646
$ws.Range("F20").Value = 0.7166393407791887
$ws.Range("G20").Value = 0.5521703372215967
$ws.Range("H20").Value = 0.6860447767167557
$ws.Range("K20").Value = 1.799487319915272
$ws.Range("L20").Value = 0.1266297481571037
$ws.Range("M20").Value = 0.5137745239966875
$ws.Range("O20").Value = 2.449609521020591
$ws.Range("C21").Value = 0.39826532030321
$ws.Range("D21").Value = 0.05071646919162731
$ws.Range("E21").Value = 0.1529691150145496
$ws.Range("F21").Value = 0.7141527124164995
$ws.Range("G21").Value = 0.5468689716337707
$ws.Range("H21").Value = 0.6780286694153119
$ws.Range("K21").Value = 2.011787271640515
$ws.Range("L21").Value = 0.1254700398979409
$ws.Range("M21").Value = 0.5595794673616155
$ws.Range("O21").Value = 2.421934138622845
$ws.Range("C22").Value = 0.4021102861833867
$ws.Range("D22").Value = 0.05300946114040528
$ws.Range("E22").Value = 0.1531627818889589
$ws.Range("F22").Value = 0.7130489578387156
$ws.Range("G22").Value = 0.5438814625187973
$ws.Range("H22").Value = 0.6731504452750983
$ws.Range("K22").Value = 2.15020656898588
$ws.Range("L22").Value = 0.124776786999572
$ws.Range("M22").Value = 0.589523873567984
$ws.Range("O22").Value = 2.405604649356803
$ws.Range("C23").Value = 0.4000412980263661
$ws.Range("D23").Value = 0.05178649276923863
$ws.Range("E23").Value = 0.1530523943652433
$ws.Range("F23").Value = 0.7135902636011693
$ws.Range("G23").Value = 0.5454321689164914
$ws.Range("H23").Value = 0.6757210098459012
$ws.Range("K23").Value = 2.076360772318765
$ws.Range("L23").Value = 0.1251408665837417
$ws.Range("M23").Value = 0.5735414087831145
$ws.Range("O23").Value = 2.414158849390304
$ws.Range("C24").Value = 0.3925759127166941
$ws.Range("D24").Value = 0.04713640149238074
$ws.Range("E24").Value = 0.1527886050046767
$ws.Range("F24").Value = 0.7166887677060245
$ws.Range("G24").Value = 0.5522650161542089
$ws.Range("H24").Value = 0.6861818806326596
$ws.Range("K24").Value = 1.796012001250631
$ws.Range("L24").Value = 0.1266497940890687
$ws.Range("M24").Value = 0.5130260453922801
$ws.Range("O24").Value = 2.450091471721748
$ws.Range("C25").Value = 0.3853164922614241
$ws.Range("D25").Value = 0.0420831207929524
$ws.Range("E25").Value = 0.1528335449171365
$ws.Range("F25").Value = 0.7222613499837109
$ws.Range("G25").Value = 0.5616725459128702
$ws.Range("H25").Value = 0.6990119689856229
$ws.Range("K25").Value = 1.492280082914419
$ws.Range("L25").Value = 0.1285547250857668
$ws.Range("M25").Value = 0.4478039359369035
$ws.Range("O25").Value = 2.496367748167728
